$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F2: ResponseMapKeys value for "Get User Account" test case
# add the new quota key to the semicolon separated list
$ws.Range("F2").Value = "account.email_address;account.is_locked;account.quotas.api_signature_requests_left;account.quotas.sms_verifications_left"

# Update G2: ResponseMapValues value for "Get User Account" test case
# add the corresponding new quota value
$ws.Range("G2").Value = "aashish.kumar@sofbang.com;false;5000;5000"

# Move the active selection to G3 (as last selected by the author before save)
$ws.Range("G3").Select()
